$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configuration")

$ws.Range("G2").Value = "API"
$ws.Range("B2").Value = "https://api-qa.metadog.racing"
